$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.027.33'
$ws.Range("E2").Value = '  +0.44%  '
$ws.Range("D3").Value = '1.561.98'
$ws.Range("E3").Value = '  +0.79%  '
$ws.Range("E4").Value = '  +0.37%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '208.33'
$ws.Range("E5").Value = '  +0.85%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.490'
$ws.Range("E6").Value = '  +0.56%  '
$ws.Range("E7").Value = '  +0.37%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.03'
$ws.Range("E8").Value = '  -0.97%  '
$ws.Range("E9").Value = '  +1.35%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0597'
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0857'
$ws.Range("E11").Value = '  +0.22%  '
$ws.Range("D12").Value = '1.786.64'
$ws.Range("E12").Value = '  +0.90%  '
$ws.Range("D13").Value = '1.584.09'
$ws.Range("E13").Value = '  +2.21%  '
$ws.Range("E14").Value = '  +0.21%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.519'
$ws.Range("E15").Value = '  +0.13%  '
$ws.Range("D16").Value = '27.026.74'
$ws.Range("E16").Value = '  +0.46%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '61.85'
$ws.Range("E17").Value = '  +0.24%  '
$ws.Range("D18").Value = '0.0₃0705'
$ws.Range("E18").Value = '  +1.00%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '215.90'
$ws.Range("E19").Value = '  -0.86%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.37'
$ws.Range("E20").Value = '  +1.71%  '
$ws.Range("E21").Value = '  +0.33%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("E22").Value = '  +2.01%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.21'
$ws.Range("E23").Value = '  -0.14%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.94'
$ws.Range("E24").Value = '  -0.52%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.58'
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("E26").Value = '  -0.34%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '15.05'
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("E28").Value = '  +1.54%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  +1.39%  '
$ws.Range("E31").Value = '  +3.92%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.22'
$ws.Range("E32").Value = '  -0.08%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.20'
$ws.Range("E33").Value = '  +4.17%  '
$ws.Range("D34").Value = '1.422.96'
$ws.Range("E34").Value = '  +0.40%  '
$ws.Range("E35").Value = '  +11.53%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.60'
$ws.Range("E36").Value = '  +1.21%  '
$ws.Range("E37").Value = '  +2.52%  '
$ws.Range("E38").Value = '  +1.46%  '
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.810'
$ws.Range("E40").Value = '  +0.23%  '
$ws.Range("E41").Value = '  +1.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  +0.29%  '
$ws.Range("E43").Value = '  +0.03%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.999'
$ws.Range("E44").Value = '  -0.34%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '64.62'
$ws.Range("E45").Value = '  +0.30%  '
$ws.Range("E46").Value = '  +0.10%  '
$ws.Range("D47").Value = '1.699.12'
$ws.Range("E47").Value = '  +0.88%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '86.55'
$ws.Range("E48").Value = '  -1.13%  '
$ws.Range("D49").Value = '0.0₆0102'
$ws.Range("E49").Value = '  +1.35%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0517'
$ws.Range("E50").Value = '  -0.01%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0960'
$ws.Range("E51").Value = '  +0.29%  '
